$d = $word.ActiveDocument

$replacements = @(
    @("97÷6=", "52÷9="),
    @("69÷3=", "96÷3="),
    @("87÷5=", "78÷3="),
    @("25÷9=", "84÷6="),
    @("38÷8=", "72÷7="),
    @("52÷2=", "48÷4="),
    @("15÷2=", "47÷4="),
    @("40÷7=", "16÷4="),
    @("77÷7=", "11÷5="),
    @("35÷4=", "31÷8="),
    @("73÷4=", "86÷8="),
    @("82÷8=", "77÷8="),
    @("41÷4=", "35÷2="),
    @("20÷2=", "44÷4="),
    @("69÷5=", "54÷9="),
    @("38÷9=", "53÷2="),
    @("95÷3=", "39÷4="),
    @("17÷2=", "70÷5="),
    @("92÷4=", "47÷2="),
    @("62÷2=", "64÷8="),
    @("58÷3=", "21÷6="),
    @("35÷6=", "24÷7="),
    @("43÷6=", "20÷3="),
    @("56÷7=", "89÷6="),
    @("99÷6=", "47÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
